$d = $word.ActiveDocument

# 1. Collapse the "Singleton design pattern for ScoreManager" paragraph
#    down to a single space, keeping the first run's (bold) formatting.
$p3 = $d.Paragraphs.Item(3)
[void]$p3.Range.Find.Execute("Singleton design pattern for ScoreManager", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2)

# 2. Move the hidden "_GoBack" bookmark from the end of the
#    "ConsoleRenderer ... PrintNewLine" paragraph to the end of the
#    "CommandParser ... 1 method" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$p7 = $d.Paragraphs.Item(7)
$endRange = $p7.Range
$endRange.Collapse(0)
[void]$endRange.MoveEnd(1, -1)
$endRange.Collapse(0)

# Work around a positioning quirk when adding a zero-length bookmark
# exactly at a paragraph's content boundary: insert a placeholder
# character, wrap the bookmark around it, then remove the placeholder
# again (the now-empty bookmark stays put, like in real Word).
$endRange.InsertAfter("X")
$markRange = $d.Range($endRange.Start, $endRange.Start + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$markRange2 = $d.Range($endRange.Start, $endRange.Start + 1)
$markRange2.Delete()
